$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 51) appended after the last existing row (row 50).
# Columns A-D are text (date/time/weekday/week stored as strings, matching
# the existing rows above). Temporarily format the date/week-number
# looking cells as Text so Excel does not auto-convert them into date or
# numeric serials, then clear the formatting again so the new row ends up
# with the same (default, unstyled) look as the rest of the data rows.
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "2024-01-12"
$ws.Range("A51").ClearFormats()

$ws.Range("B51").Value = "17:48:55"

$ws.Range("C51").Value = "Friday"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "01"
$ws.Range("D51").ClearFormats()

# Columns E-T are plain numbers.
$ws.Range("E51").Value = 135528
$ws.Range("F51").Value = 142752
$ws.Range("G51").Value = 172012
$ws.Range("H51").Value = 148321
$ws.Range("I51").Value = -1
$ws.Range("J51").Value = 119644
$ws.Range("K51").Value = 225008
$ws.Range("L51").Value = 252931
$ws.Range("M51").Value = 185015
$ws.Range("N51").Value = 110431
$ws.Range("O51").Value = 40910
$ws.Range("P51").Value = 30903
$ws.Range("Q51").Value = 73052
$ws.Range("R51").Value = -1
$ws.Range("S51").Value = 42610
$ws.Range("T51").Value = -1
